$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5 from 2023-09-06 (45175)
# to 2023-09-14 (45183), keeping the existing date formatting.
$newDate = [DateTime]::FromOADate(45183)

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
